$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set up new column C (wrap text, width matching the Polish column B style)
$ws.Columns.Item(3).ColumnWidth = 51.25

# Row by row: fix the B10 typo ("Przyczpia" -> "Przyczepia") then fill in the
# English quiz-description column C for every data row.
$ws.Range("C2").Value = "It is a flattened triangular bone with rounded edges, incorporated into the quadriceps tendon and located at the lower end of the femur."
$ws.Range("C2").WrapText = $true
$ws.Range("C3").Value = "It is a bone with a structural function that forms the lower leg. It belongs to long bones. The distal end of the tibia forms the medial malleolus."
$ws.Range("C3").WrapText = $true
$ws.Range("C4").Value = "It is the longest bone in the human body. Like any long bone, it consists of a shaft and two ends. On the back surface of the body there is a prominent rough line (linea aspera), composed of two lips: the medial one (labium mediale) and the lateral one (labium laterale). Downwards, both lips spread apart and limit the popliteal surface (facies poplitea). The medial lip ends at the bottom with the adductor tubercle (tuberculum adductorium)."
$ws.Range("C4").WrapText = $true
$ws.Range("C5").Value = "The patellar ligament connects it to the tibia. During the flexion and extension movements of the knee joint, the patella moves together with the tibia. When the knee is straight and the extensor muscles are tense, the lower part of the articular surface of the patella rests on the patella surface of the femur, the upper part lies above it; when the knee is bent, the patella moves down and back, lies in the groove between the two condyles of the femur and is immobilized. When the knee is straight and the muscles are relaxed, the patella can be moved sideways. Below the articular surface there is a convex rough area, the lower part of which serves as an attachment for the patellar ligament."
$ws.Range("C5").WrapText = $true
$ws.Range("C6").Value = "It runs from the medial surface of the lateral condyle of the femur anteromedial to the anterior intercondylar area of the tibia."
$ws.Range("C6").WrapText = $true
$ws.Range("C7").Value = "It runs from the medial surface of the medial condyle down to the posterior intercondylar area of the tibia."
$ws.Range("C7").WrapText = $true
$ws.Range("C8").Value = "The transverse knee ligament spans between the most forward points of both menisci, connecting them. It is a thin, round ligament, often inhibited in development. Sometimes it may be completely missing. It tightens when the joint rotates outwards."
$ws.Range("C8").WrapText = $true
$ws.Range("C9").Value = "Semicircular, C-shaped. Attached to the anterior and posterior intercondylar area. Firmly attached medially to the joint capsule, and laterally - to the tibial collateral ligament, which limits its mobility"
$ws.Range("C9").WrapText = $true
$ws.Range("B10").Value = "Prawie całkowicie kolista. Przyczepia się do pola międzykłykciowego przedniego i tylnego. Nie jest przytwierdzona do torebki stawowej, więc jest bardziej ruchoma od łąkotki przyśrodkowej"
$ws.Range("C10").Value = "Almost completely circular. It attaches to the anterior and posterior intercondylar areas. It is not attached to the joint capsule, so it is more mobile than the medial meniscus"
$ws.Range("C10").WrapText = $true
$ws.Range("C11").Value = "The medial condyle of the knee, also called the medial tubercle or medial process, is an important anatomical structure within the knee joint. This is a bony protrusion located on the inside of the femur, close to the knee joint. The medial condyle is crucial to the stability and function of the knee joint"
$ws.Range("C11").WrapText = $true
$ws.Range("C12").Value = "It is located on the medial side of the shin and is the longest bone of the skeleton after the femur. It has a triangular shape in cross-section. At the top, where it is involved in the formation of the knee joint, the bone thickens significantly, tapers downwards, and then widens again, although to a lesser extent than at the top. Like any long bone, it consists of a shaft and two ends"
$ws.Range("C12").WrapText = $true
$ws.Range("C13").Value = "The articular surface of the tibia is concave and adheres to the corresponding articular surface of the femur. The shape of this surface is asymmetrical, which allows for certain limitations in the mobility of the knee joint, protecting it against excessive twisting."
$ws.Range("C13").WrapText = $true
$ws.Range("C14").Value = "That is, the posterior one, above the apex, is covered with a thick layer of hyaline cartilage. It is divided into two fields, the lateral field of which is larger than the medial one, both of these parts are divided by a longitudinally running hill. It corresponds to a longitudinal groove on the patella surface of the femur."
$ws.Range("C14").WrapText = $true
$ws.Range("C15").Value = "The flat sagittal articular surface is located on the lateral condyle of the tibia. The articular surface of the fibula is directed downwards, posteriorly and laterally and is adjacent to the articular surface of the head of the fibula."
$ws.Range("C15").WrapText = $true
$ws.Range("C16").Value = "It extends from the lateral surface of the femur to the so-called the head of the fibula and is responsible for the stability of the knee from the lateral side"
$ws.Range("C16").WrapText = $true
$ws.Range("C17").Value = "It extends from the medial surface of the femur to the medial surface of the tibia, hence it is responsible for the stability of the knee joint from the medial side (i.e. from the inside)."
$ws.Range("C17").WrapText = $true
$ws.Range("C18").Value = "The ligament attaches near the posterior attachment of the lateral meniscus. It runs upwards and medially posteriorly from the posterior cruciate ligament. It usually joins it, ending at its attachment to the inner surface of the medial condyle of the femur."
$ws.Range("C18").WrapText = $true

# Row 5 grows taller once the long English quadriceps/patella description wraps
$ws.Rows.Item(5).RowHeight = 172.8

# Deselect the previously selected element after an answer in the selection quiz:
# move the active selection down to the last populated cell, C18.
$ws.Range("C18").Select()
